$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: two new columns, same bold/centered/boxed style as the
#     existing header cells ---
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "EmissionTax/EmissionType"
$ws.Range("G1").Value = "EmissionTax/EmissionTax"

# --- Row 2 (Biogas): the only fuel with emission-tax data so far ---
$ws.Range("F2").Value = "CO2"
$ws.Range("G2").Value = 24.64

# --- Rows 3-7: columns exist (widening the table) but no emission-tax
#     data has been joined in yet for these fuels, so they are blank -
#     present as an empty text value, not simply absent, matching the
#     dataframe export that produced this sheet. Plain (unstyled) cells. ---
$ws.Range("F3:G7").Formula = "'"
$ws.Range("A3").Copy()
$ws.Range("F3:G7").PasteSpecial(-4122)  # xlPasteFormats (resets style to default)
